$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add three more years of data (only the year column is populated for each new row)
$ws.Range("A7").Value = 2016
$ws.Range("A8").Value = 2015
$ws.Range("A9").Value = 2014

# Move the active selection to A9, matching the post-edit cursor position
$ws.Range("A9").Select()
